# Update the dSF column (F) values for rows 2-36 as part of a
# repull/recalculation of the mean-based statistic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -2
    3  = -1
    4  = 6
    5  = 2
    6  = -3
    7  = 1
    8  = -2
    9  = 0
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = -1
    15 = 1
    16 = 6
    17 = 2
    18 = -3
    19 = 3
    20 = -3
    21 = -4
    22 = 1
    23 = 1
    24 = 5
    25 = -1
    26 = 1
    27 = 6
    28 = 1
    29 = -2
    30 = -1
    31 = 2
    32 = 1
    33 = 1
    34 = 0
    35 = -2
    36 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
